$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 2 and row 3 course/student counts (give preference to course with more students)
$ws.Range("B2").Value = "SOC205"
$ws.Range("E2").Value = 40
$ws.Range("B3").Value = "MAT141"
$ws.Range("E3").Value = 30

# Venue reassignments
$ws.Range("F4").Value = "SLT"
$ws.Range("G4").Value = 100

$ws.Range("F5").Value = "KDLT"
$ws.Range("G5").Value = 150

$ws.Range("F6").Value = "NFLT"
$ws.Range("G6").Value = 400

$ws.Range("F7").Value = "CBN"
$ws.Range("G7").Value = 500

$ws.Range("F8").Value = "No suitable venue"
$ws.Range("G8").Value = "null"

$ws.Range("F10").Value = "SLT"
$ws.Range("G10").Value = 100

# Swap row 12 and row 13 course/student counts, and reassign venues
$ws.Range("B12").Value = "PSY101"
$ws.Range("E12").Value = 65
$ws.Range("F12").Value = "KDLT"
$ws.Range("G12").Value = 150

$ws.Range("B13").Value = "BIO101"
$ws.Range("E13").Value = 45
$ws.Range("F13").Value = "NFLT"
$ws.Range("G13").Value = 400

$ws.Range("F15").Value = "CBN"
$ws.Range("G15").Value = 500
